# Apply targeted numeric updates to Sheet1 as described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E3").Value  = 13.33
$ws.Range("E4").Value  = 12.919
$ws.Range("E7").Value  = 13.35
$ws.Range("E8").Value  = 12.913
$ws.Range("C11").Value = -12.684
$ws.Range("C12").Value = -12.745
$ws.Range("E12").Value = 12.762
$ws.Range("E14").Value = 12.911
$ws.Range("C15").Value = -12.251
$ws.Range("E22").Value = 12.862
